# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Column D ("Price") values are stored as text in the sheet (e.g. "1.0000",
# "29.582.77") even though they look numeric, so each D write is prefixed
# with a leading apostrophe (Excel's quote-prefix convention) to force
# text storage and keep exact formatting/precision instead of letting COM
# silently coerce the value to a Double and drop trailing zeros / dots.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.582.77"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").Value = "'1.854.51"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'243.71"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D6").Value = "'0.6401"
$ws.Range("E6").Value = "  -0.52%  "

$ws.Range("D7").Value = "'1.0000"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'48.53"
$ws.Range("E8").Value = "  +3.18%  "

$ws.Range("D9").Value = "'0.07555"
$ws.Range("E9").Value = "  +0.65%  "

$ws.Range("D10").Value = "'0.3006"
$ws.Range("E10").Value = "  +1.17%  "

$ws.Range("D11").Value = "'24.46"
$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("D12").Value = "'0.07668"
$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").Value = "'1.900.11"
$ws.Range("E13").Value = "  +1.96%  "

$ws.Range("D14").Value = "'5.054"
$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("D15").Value = "'0.6901"
$ws.Range("E15").Value = "  -0.15%  "

$ws.Range("D16").Value = "'84.10"
$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("D17").Value = "'0.000009634"
$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("D18").Value = "'2.156.32"
$ws.Range("E18").Value = "  +1.92%  "

$ws.Range("D19").Value = "'6.275"
$ws.Range("E19").Value = "  +3.28%  "

$ws.Range("D20").Value = "'29.636.72"
$ws.Range("E20").Value = "  -0.70%  "

$ws.Range("D21").Value = "'238.90"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").Value = "'12.64"
$ws.Range("E22").Value = "  -0.38%  "

$ws.Range("D23").Value = "'0.9996"

$ws.Range("D24").Value = "'7.668"
$ws.Range("E24").Value = "  +3.00%  "

$ws.Range("D25").Value = "'1.001"
$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("D26").Value = "'157.26"
$ws.Range("E26").Value = "  -0.93%  "

$ws.Range("D27").Value = "'0.1405"
$ws.Range("E27").Value = "  -1.95%  "

$ws.Range("D28").Value = "'8.511"
$ws.Range("E28").Value = "  -0.28%  "

$ws.Range("D29").Value = "'17.84"
$ws.Range("E29").Value = "  -0.63%  "

$ws.Range("D30").Value = "'1.489"
$ws.Range("E30").Value = "  -0.37%  "

$ws.Range("D31").Value = "'0.05908"
$ws.Range("E31").Value = "  -3.73%  "

$ws.Range("D32").Value = "'1.284"
$ws.Range("E32").Value = "  +1.00%  "

$ws.Range("D33").Value = "'4.149"
$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("D34").Value = "'4.087"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("D35").Value = "'1.931"
$ws.Range("E35").Value = "  +2.74%  "

$ws.Range("D36").Value = "'1.183"
$ws.Range("E36").Value = "  +1.04%  "

$ws.Range("D37").Value = "'0.7273"
$ws.Range("E37").Value = "  -0.99%  "

$ws.Range("E38").Value = "  -0.48%  "

$ws.Range("D39").Value = "'2.802"
$ws.Range("E39").Value = "  -1.72%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01779"
$ws.Range("E40").Value = "  -0.87%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "'1.213.22"
$ws.Range("E41").Value = "  -0.20%  "

$ws.Range("D42").Value = "'0.9164"
$ws.Range("E42").Value = "  -1.18%  "

$ws.Range("D43").Value = "'6.125"
$ws.Range("E43").Value = "  -1.01%  "

$ws.Range("D44").Value = "'2.061.75"
$ws.Range("E44").Value = "  +1.89%  "

$ws.Range("D45").Value = "'0.9997"
$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("D46").Value = "'102.07"
$ws.Range("E46").Value = "  -0.19%  "

$ws.Range("D47").Value = "'67.44"
$ws.Range("E47").Value = "  +1.69%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.500"
$ws.Range("E48").Value = "  +11.65%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000122"
$ws.Range("E49").Value = "  +0.47%  "

$ws.Range("D50").Value = "'0.4072"
$ws.Range("E50").Value = "  -0.19%  "

$ws.Range("D51").Value = "'9.188"
$ws.Range("E51").Value = "  -0.61%  "

